$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '40.954.03'
$ws.Range('E2').Value = '  -2.45%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.168.10'
$ws.Range('E3').Value = '  -2.13%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '236.08'
$ws.Range('E5').Value = '  -2.92%  '
$ws.Range('E6').Value = '  -2.66%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '69.31'
$ws.Range('E7').Value = '  -5.81%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.570'
$ws.Range('E9').Value = '  -6.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.64'
$ws.Range('E10').Value = '  -8.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0920'
$ws.Range('E11').Value = '  -4.00%  '
$ws.Range('E12').Value = '  -2.76%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.71'
$ws.Range('E13').Value = '  -5.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.494.25'
$ws.Range('E14').Value = '  -2.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '13.77'
$ws.Range('E15').Value = '  -3.36%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.190.00'
$ws.Range('E16').Value = '  -1.72%  '
$ws.Range('B17').Value = 'Polygon'
$ws.Range('C17').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.803'
$ws.Range('E17').Value = '  -4.64%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '40.858.32'
$ws.Range('E18').Value = '  -2.46%  '
$ws.Range('E19').Value = '  -7.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.17'
$ws.Range('E20').Value = '  -2.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.91'
$ws.Range('E21').Value = '  -4.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.63'
$ws.Range('E22').Value = '  -5.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '224.44'
$ws.Range('E23').Value = '  -2.07%  '
$ws.Range('E24').Value = '  -8.26%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.87'
$ws.Range('E26').Value = '  -6.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.52'
$ws.Range('E27').Value = '  -1.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.19'
$ws.Range('E28').Value = '  -3.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '165.78'
$ws.Range('E30').Value = '  -0.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.77'
$ws.Range('E31').Value = '  -3.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '30.78'
$ws.Range('E32').Value = '  +4.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0767'
$ws.Range('E33').Value = '  -3.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.10'
$ws.Range('E34').Value = '  -8.89%  '
$ws.Range('E35').Value = '  -3.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.103'
$ws.Range('E36').Value = '  -9.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.12'
$ws.Range('E37').Value = '  -4.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0284'
$ws.Range('E38').Value = '  -5.56%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.27'
$ws.Range('E39').Value = '  -6.00%  '
$ws.Range('E40').Value = '  -3.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.42'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '59.59'
$ws.Range('E42').Value = '  -7.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.188'
$ws.Range('E43').Value = '  -6.11%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0970'
$ws.Range('E44').Value = '  -3.95%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.23'
$ws.Range('E45').Value = '  -5.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '98.22'
$ws.Range('E46').Value = '  -6.06%  '
$ws.Range('E47').Value = '  -3.28%  '
$ws.Range('E48').Value = '  -3.74%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.21'
$ws.Range('E49').Value = '  -7.83%  '
$ws.Range('E50').Value = '  -2.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.372.79'
